# Regenerate the handback/handoff report: refresh the "Latest Handoff"
# timestamp for every file row that has actually been (re-)handed off
# in this run (i.e. every row except the one still "In Translation").
#
# Overview sheet (sheet1) keeps the date in its own "yy-dd-mm" style text,
# while the per-locale sheets (zh-cn / de-de) keep the real
# "Latest Handoff Datetime" value. Each locale's timestamp is distinct.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 11, 12, 13, 14, 15, 16)

$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("D$r").Value = "2016-23-17 22:23:35"
}

$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "2016-03-17 22:23:31"
}

$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("E$r").Value = "2016-03-17 22:23:35"
}
